$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header values for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Numeric data for column I (rows 2-17)
$iValues = @(4, 6, 9, 6, 3, 7, 7, 7, 6, 3, 8, 7, 9, 4, 9, 9)
# Numeric data for column J (rows 2-17)
$jValues = @(5, 7, 9, 6, 4, 8, 8, 8, 7, 4, 8, 7, 9, 6, 9, 9)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
